$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "cd_vendedor"
$ws.Range("Z1").Value = "total"

$ws.Range("Z2:Z12").Formula = "=SUM(B2:Y2)"
